$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.0202168574806743
$ws.Range("E4").Value = 0.009442948505717938
$ws.Range("F4").Value = 0.01953356698934268
$ws.Range("H4").Value = -0.01932747370109895
$ws.Range("J4").Value = -0.006047168476769442
$ws.Range("C5").Value = -0.009474448218977928
$ws.Range("E5").Value = -0.006712346380493854
$ws.Range("F5").Value = 0.002650578058023122
$ws.Range("H5").Value = 0.001057213386288535
$ws.Range("J5").Value = 0.002319020709093875
$ws.Range("C6").Value = 0.005435377465415098
$ws.Range("E6").Value = -0.0004470717298828692
$ws.Range("F6").Value = 0.003988637151545486
$ws.Range("H6").Value = -0.004334404205376167
$ws.Range("J6").Value = -0.001797001047756817
$ws.Range("C7").Value = 0.003723370708934828
$ws.Range("E7").Value = -0.02004915440196617
$ws.Range("F7").Value = 0.0005655055906202235
$ws.Range("H7").Value = 0.005264370930574836
$ws.Range("J7").Value = 0.0007928998079122779
$ws.Range("C8").Value = 0.078027200401088
$ws.Range("E8").Value = -0.2246937710037508
$ws.Range("F8").Value = -0.9628999603399983
$ws.Range("H8").Value = 0.9999999085119963
$ws.Range("J8").Value = -0.05410777799973252
$ws.Range("C9").Value = 0.9709073071242921
$ws.Range("E9").Value = 0.01628250861930034
$ws.Range("F9").Value = 0.02394311058972442
$ws.Range("H9").Value = -0.02724026144161045
$ws.Range("J9").Value = 0.01998509579277222
$ws.Range("C10").Value = 0.005281874131274964
$ws.Range("E10").Value = -0.01527151433886057
$ws.Range("F10").Value = -0.006659213738368548
$ws.Range("H10").Value = 0.01121673366466935
$ws.Range("J10").Value = -0.02628444167791655
$ws.Range("C11").Value = -0.002057418226296729
$ws.Range("E11").Value = 0.008232176009287039
$ws.Range("F11").Value = -0.007694493427779736
$ws.Range("H11").Value = 0.005147597965903918
$ws.Range("J11").Value = 0.0003062826501658919
$ws.Range("C12").Value = 0.06920739944029597
$ws.Range("E12").Value = 0.0001130480685219227
$ws.Range("F12").Value = -0.06565117481804698
$ws.Range("H12").Value = 0.02383479896939196
$ws.Range("J12").Value = -0.006597131294680632
$ws.Range("C13").Value = 0.06714615849384634
$ws.Range("E13").Value = 0.851684828627393
$ws.Range("F13").Value = -0.2456314848332594
$ws.Range("H13").Value = -0.01319485435179417
$ws.Range("J13").Value = 0.007292730773770639
$ws.Range("C14").Value = -0.163567675598707
$ws.Range("E14").Value = 0.02252694186107767
$ws.Range("F14").Value = -0.02517815678312627
$ws.Range("H14").Value = 0.0207202950368118
$ws.Range("J14").Value = 0.002127535125176377
$ws.Range("C15").Value = -0.01453517866140715
$ws.Range("E15").Value = 0.008239290665571625
$ws.Range("F15").Value = -0.003166207422648297
$ws.Range("H15").Value = 0.002036556657462266
$ws.Range("J15").Value = -0.008950762962789171
$ws.Range("C16").Value = -0.005310071540402861
$ws.Range("E16").Value = -0.02103732688949307
$ws.Range("F16").Value = -0.02111880468475218
$ws.Range("H16").Value = 0.02683267364930694
$ws.Range("J16").Value = 0.01543438316297372
$ws.Range("C17").Value = 0.006002473200098927
$ws.Range("E17").Value = 0.01879000193560008
$ws.Range("F17").Value = -0.04778786303151451
$ws.Range("H17").Value = 0.04119389041575561
$ws.Range("J17").Value = 0.01000934417930227
$ws.Range("C18").Value = 0.02663205562528222
$ws.Range("E18").Value = 0.01569075768363031
$ws.Range("F18").Value = -0.01690509658020386
$ws.Range("H18").Value = 0.01337378808695152
$ws.Range("J18").Value = 0.005914755413071732
$ws.Range("C19").Value = 0.01085104728204189
$ws.Range("E19").Value = 0.0005848441193937647
$ws.Range("F19").Value = -0.003659036690361467
$ws.Range("H19").Value = 0.002064199858567994
$ws.Range("J19").Value = -0.009163207042856778
$ws.Range("C20").Value = 0.009232352721294109
$ws.Range("E20").Value = -0.005457511322300452
$ws.Range("F20").Value = 0.001656558306262332
$ws.Range("H20").Value = -0.000583583351343334
$ws.Range("J20").Value = 0.01209684775406746
$ws.Range("C21").Value = 0.0236794775871791
$ws.Range("E21").Value = 0.006258251578330062
$ws.Range("F21").Value = -0.02467307561092302
$ws.Range("H21").Value = 0.02410193568407742
$ws.Range("J21").Value = -0.005996936327910488
$ws.Range("C22").Value = 0.01217704704708188
$ws.Range("E22").Value = -0.001811484840459393
$ws.Range("F22").Value = 0.001723300964932038
$ws.Range("H22").Value = -0.0007693279987731199
$ws.Range("J22").Value = 0.0260608936875783
$ws.Range("C23").Value = -0.008212617160504684
$ws.Range("E23").Value = 0.002784168975366759
$ws.Range("F23").Value = 0.009293538419741535
$ws.Range("H23").Value = -0.008494716339788653
$ws.Range("J23").Value = 0.01308304582140079
